$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5:D102").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# New column D values (most recent reporting period)
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 465800
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -25900
$ws.Range("D17").Value = 58800
$ws.Range("D18").Value = 407000
$ws.Range("D20").Value = -134500
$ws.Range("D21").Value = 298500
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 272600
$ws.Range("D24").Value = 58200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 214400
$ws.Range("D27").Value = 214400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 1600
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 134500
$ws.Range("D33").Value = 215900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 215900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 316800
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 506900
$ws.Range("D49").Value = 282500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 11872000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 160400
$ws.Range("D62").Value = 4800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 9932400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2064100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1939600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 215900
$ws.Range("D83").Value = 25900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 229800
$ws.Range("D91").Value = -21400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 379400
$ws.Range("D96").Value = -49600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -557800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 51400

# Corrections to shifted values that differ from a pure shift
$ws.Range("E62").Value = 5200
$ws.Range("F62").Value = 1600
$ws.Range("J91").Value = -62000
